$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "60.428.30"
$ws.Range("E2").Value = "  +2.53%  "

$ws.Range("D3").Value = "2.685.37"
$ws.Range("E3").Value = "  +1.32%  "

$ws.Range("E4").Value = "  -0.06%  "

Set-TextValue $ws.Range("D5") "522.83"
$ws.Range("E5").Value = "  +1.55%  "

Set-TextValue $ws.Range("D6") "146.18"
$ws.Range("E6").Value = "  +1.31%  "

Set-TextValue $ws.Range("D7") "0.997"
$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("E8").Value = "  +1.35%  "

$ws.Range("D9").Value = "2.706.38"
$ws.Range("E9").Value = "  +0.96%  "

$ws.Range("E10").Value = "  +2.40%  "

$ws.Range("E11").Value = "  -0.11%  "

$ws.Range("E12").Value = "  +1.75%  "

$ws.Range("E13").Value = "  +1.61%  "

$ws.Range("D14").Value = "3.157.95"
$ws.Range("E14").Value = "  +1.36%  "

$ws.Range("D15").Value = "60.444.28"
$ws.Range("E15").Value = "  +2.62%  "

$ws.Range("E16").Value = "  +0.71%  "

$ws.Range("D17").Value = "2.765.30"
$ws.Range("E17").Value = "  +3.42%  "

$ws.Range("E18").Value = "  +1.26%  "

Set-TextValue $ws.Range("D19") "351.01"
$ws.Range("E19").Value = "  +2.96%  "

$ws.Range("E20").Value = "  -0.07%  "

Set-TextValue $ws.Range("D21") "10.61"
$ws.Range("E21").Value = "  +1.53%  "

Set-TextValue $ws.Range("D22") "6.32"
$ws.Range("E22").Value = "  +3.03%  "

$ws.Range("E23").Value = "  +0.07%  "

Set-TextValue $ws.Range("D24") "62.90"
$ws.Range("E24").Value = "  +2.99%  "

$ws.Range("E25").Value = "  +0.42%  "

Set-TextValue $ws.Range("D26") "0.169"
$ws.Range("E26").Value = "  +5.09%  "

Set-TextValue $ws.Range("D27") "0.994"
$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("E28").Value = "  +0.48%  "

Set-TextValue $ws.Range("D29") "7.28"
$ws.Range("E29").Value = "  +1.78%  "

$ws.Range("E30").Value = "  +7.70%  "

$ws.Range("E31").Value = "  +0.12%  "

$ws.Range("E33").Value = "  +0.74%  "

Set-TextValue $ws.Range("D34") "148.44"
$ws.Range("E34").Value = "  -0.51%  "

Set-TextValue $ws.Range("D35") "4.33"
$ws.Range("E35").Value = "  +7.26%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D36") "1.23"
$ws.Range("E36").Value = "  +7.32%  "

$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D37") "0.949"
$ws.Range("E37").Value = "  -7.22%  "

$ws.Range("E38").Value = "  +10.19%  "

$ws.Range("E39").Value = "  +2.49%  "

Set-TextValue $ws.Range("D40") "36.86"
$ws.Range("E40").Value = "  +0.67%  "

Set-TextValue $ws.Range("D41") "3.69"
$ws.Range("E41").Value = "  +0.32%  "

Set-TextValue $ws.Range("D42") "282.29"
$ws.Range("E42").Value = "  -0.32%  "

$ws.Range("E43").Value = "  -1.05%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D44") "0.996"
$ws.Range("E44").Value = "  +0.34%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D45") "19.95"
$ws.Range("E45").Value = "  +2.10%  "

$ws.Range("E46").Value = "  +0.56%  "

$ws.Range("D47").Value = "2.126.56"
$ws.Range("E47").Value = "  +6.75%  "

$ws.Range("E48").Value = "  +0.96%  "

$ws.Range("E49").Value = "  +3.33%  "

$ws.Range("E50").Value = "  +2.01%  "

Set-TextValue $ws.Range("D51") "10.46"
$ws.Range("E51").Value = "  +1.78%  "
